# Updates cryptos list data (Coin/Link/Price/Volume columns) per
# "Updated cryptos list on Sat May 13 14:41:27 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.431.57'
$ws.Range("E2").Value = '  +3.52%  '
$ws.Range("D3").Value = '1.840.71'
$ws.Range("E3").Value = '  +3.57%  '
$ws.Range("D4").Value = '1.028'
$ws.Range("E4").Value = '  +2.71%  '
$ws.Range("D5").Value = "'318.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.53%  '
$ws.Range("D6").Value = '1.025'
$ws.Range("E6").Value = '  +2.41%  '
$ws.Range("D7").Value = '0.4363'
$ws.Range("E7").Value = '  +1.74%  '
$ws.Range("D8").Value = '0.3731'
$ws.Range("E8").Value = '  +2.63%  '
$ws.Range("D9").Value = '0.07356'
$ws.Range("D10").Value = '0.8764'
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").Value = '21.42'
$ws.Range("E11").Value = '  +4.04%  '
$ws.Range("D12").Value = '1.979.14'
$ws.Range("E12").Value = '  +11.93%  '
$ws.Range("D13").Value = '5.484'
$ws.Range("E13").Value = '  +4.06%  '
$ws.Range("D14").Value = '6.686'
$ws.Range("E14").Value = '  +3.54%  '
$ws.Range("D15").Value = '0.07171'
$ws.Range("E15").Value = '  +4.31%  '
$ws.Range("D16").Value = '82.29'
$ws.Range("E16").Value = '  +4.11%  '
$ws.Range("D17").Value = '1.029'
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").Value = '0.000009001'
$ws.Range("D19").Value = '1.025'
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").Value = '15.41'
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").Value = '27.456.66'
$ws.Range("D22").Value = '5.257'
$ws.Range("E22").Value = '  +2.60%  '
$ws.Range("D23").Value = '11.16'
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("D24").Value = '2.177.27'
$ws.Range("E24").Value = '  +8.57%  '
$ws.Range("D25").Value = '157.05'
$ws.Range("E25").Value = '  +3.17%  '
$ws.Range("D26").Value = '1.901'
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("D27").Value = '18.56'
$ws.Range("E27").Value = '  +2.74%  '
$ws.Range("D28").Value = '5.277'
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("D29").Value = '1.923'
$ws.Range("E29").Value = '  +6.89%  '
$ws.Range("D30").Value = '115.47'
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("D31").Value = '0.09032'
$ws.Range("E31").Value = '  +1.36%  '
$ws.Range("D32").Value = '1.201'
$ws.Range("E32").Value = '  +6.27%  '
$ws.Range("D33").Value = '0.7594'
$ws.Range("E33").Value = '  +4.03%  '
$ws.Range("D34").Value = '4.475'
$ws.Range("E34").Value = '  +3.25%  '
$ws.Range("D35").Value = "'2.860"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.03%  '
$ws.Range("D36").Value = '1.027'
$ws.Range("E36").Value = '  +2.70%  '
$ws.Range("E37").Value = '  +4.66%  '
$ws.Range("D38").Value = '0.01959'
$ws.Range("E38").Value = '  +3.51%  '
$ws.Range("D39").Value = '0.05255'
$ws.Range("E39").Value = '  +1.76%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5171'
$ws.Range("E40").Value = '  +4.58%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.792'
$ws.Range("E41").Value = '  +5.26%  '
$ws.Range("D42").Value = '0.1662'
$ws.Range("E42").Value = '  +2.81%  '
$ws.Range("D43").Value = '6.538'
$ws.Range("E43").Value = '  +3.28%  '
$ws.Range("D44").Value = '8.474'
$ws.Range("E44").Value = '  +5.52%  '
$ws.Range("D45").Value = '108.74'
$ws.Range("E45").Value = '  +3.34%  '
$ws.Range("E46").Value = '  +3.30%  '
$ws.Range("E47").Value = '  +2.68%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = "'0.4640"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.55%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.673'
$ws.Range("E49").Value = '  +1.77%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06293'
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = "'1.880"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.51%  '
